# Auto-generated Excel COM-interop script to apply profit/price updates
# across multiple worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 77272.69500000001
$ws.Range("I53").Value = 254.4
$ws.Range("J53").Value = 334000.34
$ws.Range("K53").Value = 254.4
$ws.Range("L53").Value = 334000.34
$ws.Range("M53").Value = 382.6
$ws.Range("N53").Value = -335274.34
$ws.Range("H97").Value = 949.5
$ws.Range("J97").Value = 949.5
$ws.Range("L97").Value = 2848.5
$ws.Range("N97").Value = -3840.5
$ws.Range("H108").Value = 79999
$ws.Range("I108").Value = 79999
$ws.Range("J108").Value = 0
$ws.Range("K108").Value = 79999
$ws.Range("L108").Value = 0
$ws.Range("M108").Value = -76159
$ws.Range("N108").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5648.022
$ws.Range("I32").Value = 4972.6816
$ws.Range("J32").Value = 20505.5
$ws.Range("K32").Value = 4972.6816
$ws.Range("L32").Value = 20505.5
$ws.Range("M32").Value = -4685.6816
$ws.Range("N32").Value = -21079.5
$ws.Range("H45").Value = 2395.875
$ws.Range("J45").Value = 2473.25
$ws.Range("L45").Value = 2473.25
$ws.Range("N45").Value = -3227.25
$ws.Range("H61").Value = 2561.9412
$ws.Range("I61").Value = 1465.6923
$ws.Range("J61").Value = 6124.75
$ws.Range("K61").Value = 1465.6923
$ws.Range("L61").Value = 6124.75
$ws.Range("M61").Value = -1253.6923
$ws.Range("N61").Value = -6548.75
$ws.Range("H136").Value = 2561.9412
$ws.Range("I136").Value = 1465.6923
$ws.Range("J136").Value = 6124.75
$ws.Range("K136").Value = 4397.0769
$ws.Range("L136").Value = 18374.25
$ws.Range("M136").Value = -1847.0769
$ws.Range("N136").Value = -23474.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1753.7142
$ws.Range("I86").Value = 1407.8077
$ws.Range("J86").Value = 2753
$ws.Range("K86").Value = 1407.8077
$ws.Range("L86").Value = 2753
$ws.Range("M86").Value = -284.8077000000001
$ws.Range("N86").Value = -4999
$ws.Range("H89").Value = 1753.7142
$ws.Range("I89").Value = 1407.8077
$ws.Range("J89").Value = 2753
$ws.Range("K89").Value = 7039.038500000001
$ws.Range("L89").Value = 13765
$ws.Range("M89").Value = -1423.038500000001
$ws.Range("N89").Value = -24997
$ws.Range("H94").Value = 8621658
$ws.Range("I94").Value = 8929561
$ws.Range("J94").Value = 367
$ws.Range("K94").Value = 8929561
$ws.Range("L94").Value = 367
$ws.Range("M94").Value = -8929110
$ws.Range("N94").Value = -1269
$ws.Range("H99").Value = 1608.5
$ws.Range("I99").Value = 999
$ws.Range("K99").Value = 999
$ws.Range("M99").Value = 499

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 2874.75
$ws.Range("I3").Value = 750
$ws.Range("J3").Value = 4999.5
$ws.Range("K3").Value = 750
$ws.Range("L3").Value = 4999.5
$ws.Range("M3").Value = -637
$ws.Range("N3").Value = -5225.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2537.75
$ws.Range("J5").Value = 4157.222
$ws.Range("L5").Value = 12471.666
$ws.Range("N5").Value = -12695.666
$ws.Range("H11").Value = 30081.172
$ws.Range("I11").Value = 41784.84
$ws.Range("K11").Value = 125354.52
$ws.Range("M11").Value = -125214.52
$ws.Range("H80").Value = 5498.5
$ws.Range("I80").Value = 4002
$ws.Range("J80").Value = 5997.3335
$ws.Range("K80").Value = 12006
$ws.Range("L80").Value = 17992.0005
$ws.Range("M80").Value = -11070
$ws.Range("N80").Value = -19864.0005
$ws.Range("H83").Value = 5498.5
$ws.Range("I83").Value = 4002
$ws.Range("J83").Value = 5997.3335
$ws.Range("K83").Value = 36018
$ws.Range("L83").Value = 53976.0015
$ws.Range("M83").Value = -31338
$ws.Range("N83").Value = -63336.0015
$ws.Range("H114").Value = 4009.5715
$ws.Range("I114").Value = 436.2857
$ws.Range("J114").Value = 7582.857
$ws.Range("K114").Value = 1308.8571
$ws.Range("L114").Value = 22748.571
$ws.Range("M114").Value = 1945.1429
$ws.Range("N114").Value = -29256.571
$ws.Range("H135").Value = 2537.75
$ws.Range("J135").Value = 4157.222
$ws.Range("L135").Value = 37414.998
$ws.Range("N135").Value = -42484.998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()
$ws.Range("H24").Value = 16833.334
$ws.Range("J24").Value = 16833.334
$ws.Range("L24").Value = 16833.334
$ws.Range("N24").Value = -17179.334
$ws.Range("H70").Value = 6724.7427
$ws.Range("I70").Value = 5646.0415
$ws.Range("K70").Value = 5646.0415
$ws.Range("M70").Value = -5376.0415
$ws.Range("H73").Value = 6724.7427
$ws.Range("I73").Value = 5646.0415
$ws.Range("K73").Value = 5646.0415
$ws.Range("M73").Value = -4710.0415
$ws.Range("H80").Value = 3327.739
$ws.Range("J80").Value = 4940
$ws.Range("L80").Value = 4940
$ws.Range("N80").Value = -6936
$ws.Range("H83").Value = 3327.739
$ws.Range("J83").Value = 4940
$ws.Range("L83").Value = 24700
$ws.Range("N83").Value = -34684
$ws.Range("H122").Value = 2108.32
$ws.Range("I122").Value = 1774.2632
$ws.Range("K122").Value = 5322.7896
$ws.Range("M122").Value = -2872.7896
$ws.Range("H132").Value = 2762
$ws.Range("I132").Value = 2709.653
$ws.Range("K132").Value = 8128.958999999999
$ws.Range("M132").Value = -5598.958999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 834.3
$ws.Range("I46").Value = 408.5
$ws.Range("J46").Value = 1473
$ws.Range("K46").Value = 408.5
$ws.Range("L46").Value = 1473
$ws.Range("M46").Value = -220.5
$ws.Range("N46").Value = -1849
$ws.Range("H68").Value = 4908.643
$ws.Range("I68").Value = 3357.75
$ws.Range("K68").Value = 3357.75
$ws.Range("M68").Value = -2608.75
$ws.Range("H71").Value = 4908.643
$ws.Range("I71").Value = 3357.75
$ws.Range("K71").Value = 16788.75
$ws.Range("M71").Value = -13044.75
$ws.Range("H74").Value = 50975
$ws.Range("I74").Value = 50975
$ws.Range("K74").Value = 50975
$ws.Range("M74").Value = -49977
$ws.Range("H77").Value = 50975
$ws.Range("I77").Value = 50975
$ws.Range("K77").Value = 152925
$ws.Range("M77").Value = -147933
$ws.Range("H82").Value = 835.5454999999999
$ws.Range("I82").Value = 570.75
$ws.Range("J82").Value = 1153.3
$ws.Range("K82").Value = 570.75
$ws.Range("L82").Value = 1153.3
$ws.Range("M82").Value = -209.75
$ws.Range("N82").Value = -1875.3
$ws.Range("H85").Value = 835.5454999999999
$ws.Range("I85").Value = 570.75
$ws.Range("J85").Value = 1153.3
$ws.Range("K85").Value = 570.75
$ws.Range("L85").Value = 1153.3
$ws.Range("M85").Value = 677.25
$ws.Range("N85").Value = -3649.3
$ws.Range("H88").Value = 40189
$ws.Range("J88").Value = 40189
$ws.Range("L88").Value = 40189
$ws.Range("N88").Value = -41045
$ws.Range("H91").Value = 40189
$ws.Range("J91").Value = 40189
$ws.Range("L91").Value = 40189
$ws.Range("N91").Value = -43153
$ws.Range("H93").Value = 15874890
$ws.Range("I93").Value = 17545614
$ws.Range("J93").Value = 3004
$ws.Range("K93").Value = 17545614
$ws.Range("L93").Value = 3004
$ws.Range("M93").Value = -17544366
$ws.Range("N93").Value = -5500

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").ClearContents()
$ws.Range("H22").Value = 3000
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("H107").Value = 1346.2941
$ws.Range("I107").Value = 1048.9
$ws.Range("J107").Value = 1771.1428
$ws.Range("K107").Value = 3146.7
$ws.Range("L107").Value = 5313.428400000001
$ws.Range("M107").Value = -1226.7
$ws.Range("N107").Value = -9153.428400000001
$ws.Range("H122").Value = 1658.6888
$ws.Range("I122").Value = 1526.2069
$ws.Range("J122").Value = 1898.8125
$ws.Range("K122").Value = 4578.620699999999
$ws.Range("L122").Value = 5696.4375
$ws.Range("M122").Value = -2128.620699999999
$ws.Range("N122").Value = -10596.4375
